# RegionDF_countinuously.xlsx — "NEW EXCELS 16/4 and bug fix @ regions"
#
# 1) Bug fix: four Fans counts on 2019-04-15 (rows 268-271) were re-pulled
#    and corrected.
# 2) The last two days' worth of rows (2019-04-16 and 2019-04-17, rows
#    275-288) are removed — those days hadn't actually finished fetching
#    yet, so they get dropped until the next continuous run.
# 3) Column A is widened by one unit to keep fitting the date column, and
#    the view is scrolled back up a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Corrected Fans values for 2019-04-15 -----------------------------
$ws.Range("C268").Value = 116
$ws.Range("C269").Value = 109
$ws.Range("C270").Value = 26
$ws.Range("C271").Value = 21

# --- 2) Drop the trailing, not-yet-complete rows 275:288 ------------------
$ws.Rows("275:288").Delete()

# --- 3) Small cosmetic fixes ----------------------------------------------
# Column A ("Date Fetched") grows from raw width 10 to raw width 11.
$ws.Columns("A:A").ColumnWidth = 10.1666666667

# Scroll the view back up so row 229 is at the top again.
[void]$ws.Range("A229").Select()
